# Performance.xlsx update
#  - added array pool / use of array in moves generator
#  - refactored performance test class, new moves cache, vector vs array pool bench
#
# This reproduces, via Excel COM automation, the same content change that the
# commit made to docs/Performance.xlsx:
#   * two new benchmark "note" strings
#   * two stray numbers removed from column R (60/61)
#   * P63 note cell un-bolded (now matches the other note cells)
#   * three brand-new benchmark rows (67-69) appended on
#     "Initial Position Single Thread", with the same formula pattern as the
#     previous block of rows (59-65)
#   * viewport/selection updated to the new bottom of the sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Initial Position Single Thread")
$ws2 = $wb.Worksheets.Item("vs other Enignes")

# ---------------------------------------------------------------------------
# 1. Column R leftover numbers are gone - the two cells become blank (same
#    plain/general look as the rest of that block, e.g. C60).
# ---------------------------------------------------------------------------
$ws1.Range("R60").ClearContents()
$ws1.Range("R61").ClearContents()
$ws1.Range("R60").Style = $ws1.Range("C60").Style
$ws1.Range("R61").Style = $ws1.Range("C61").Style

# ---------------------------------------------------------------------------
# 2. P63 note cell loses its bold/variant styling - now matches the plain
#    note style used by the sibling cells (P1, P59 ...).
# ---------------------------------------------------------------------------
$ws1.Range("P63").Style = $ws1.Range("P59").Style

# ---------------------------------------------------------------------------
# 3. Three new rows of benchmark data, continuing the existing 4/5/6-depth
#    pattern (rows 59-61 and 63-65) for a new commit dated 2025-07-14
#    (serial 45862).
# ---------------------------------------------------------------------------

# --- row 67 (depth 4, mirrors row 59 / row 63) ---
$ws1.Range("A67").Style = $ws1.Range("A63").Style
$ws1.Range("B67").Style = $ws1.Range("B63").Style
$ws1.Range("C67").Style = $ws1.Range("C63").Style
$ws1.Range("D67").Style = $ws1.Range("D63").Style
$ws1.Range("E67").Style = $ws1.Range("E63").Style
$ws1.Range("F67").Style = $ws1.Range("F63").Style
$ws1.Range("G67").Style = $ws1.Range("G63").Style
$ws1.Range("H67").Style = $ws1.Range("H63").Style
$ws1.Range("I67").Style = $ws1.Range("I63").Style
$ws1.Range("J67").Style = $ws1.Range("J63").Style
$ws1.Range("K67").Style = $ws1.Range("K63").Style
$ws1.Range("L67").Style = $ws1.Range("L63").Style
$ws1.Range("M67").Style = $ws1.Range("M63").Style
$ws1.Range("N67").Style = $ws1.Range("N63").Style

$ws1.Range("A67").Value = 45862
$ws1.Range("C67").Value = 4
$ws1.Range("D67").Value = 206603
$ws1.Range("E67").Value = 274
$ws1.Range("F67").Formula = "=D67/E67*1000"
$ws1.Range("G67").Formula = "=(E63-E67)/E63"
$ws1.Range("H67").Formula = "=(F67-80000000)/80000000"
$ws1.Range("I67").Value = 4
$ws1.Range("J67").Value = 197281
$ws1.Range("K67").Value = 9
$ws1.Range("L67").Formula = "=J67/K67*1000"
$ws1.Range("M67").Formula = "=(K63-K67)/K63"
$ws1.Range("N67").Formula = "=(L67-80000000)/80000000"
$ws1.Range("P67").Value = "replaced pool vector with pool array (no board and opposite caching)"

# --- row 68 (depth 5, mirrors row 60 / row 64) ---
$ws1.Range("C68").Style = $ws1.Range("C64").Style
$ws1.Range("D68").Style = $ws1.Range("D64").Style
$ws1.Range("E68").Style = $ws1.Range("E64").Style
$ws1.Range("F68").Style = $ws1.Range("F64").Style
$ws1.Range("G68").Style = $ws1.Range("G64").Style
$ws1.Range("H68").Style = $ws1.Range("H64").Style
$ws1.Range("I68").Style = $ws1.Range("I64").Style
$ws1.Range("J68").Style = $ws1.Range("J64").Style
$ws1.Range("K68").Style = $ws1.Range("K64").Style
$ws1.Range("L68").Style = $ws1.Range("L64").Style
$ws1.Range("M68").Style = $ws1.Range("M64").Style
$ws1.Range("N68").Style = $ws1.Range("N64").Style

$ws1.Range("C68").Value = 5
$ws1.Range("D68").Value = 5072212
$ws1.Range("E68").Value = 7220
$ws1.Range("F68").Formula = "=D68/E68*1000"
$ws1.Range("G68").Formula = "=(E64-E68)/E64"
$ws1.Range("H68").Formula = "=(F68-80000000)/80000000"
$ws1.Range("I68").Value = 5
$ws1.Range("J68").Value = 4880523
$ws1.Range("K68").Value = 223
$ws1.Range("L68").Formula = "=J68/K68*1000"
$ws1.Range("M68").Formula = "=(K64-K68)/K64"
$ws1.Range("N68").Formula = "=(L68-80000000)/80000000"

# --- row 69 (depth 6, mirrors row 61 / row 65) ---
$ws1.Range("I69").Style = $ws1.Range("I65").Style
$ws1.Range("J69").Style = $ws1.Range("J65").Style
$ws1.Range("K69").Style = $ws1.Range("K65").Style
$ws1.Range("L69").Style = $ws1.Range("L65").Style
$ws1.Range("M69").Style = $ws1.Range("M65").Style
$ws1.Range("N69").Style = $ws1.Range("N65").Style

$ws1.Range("I69").Value = 6
$ws1.Range("J69").Value = 119060324
$ws1.Range("K69").Value = 5953
$ws1.Range("L69").Formula = "=J69/K69*1000"
$ws1.Range("M69").Formula = "=(K65-K69)/K65"
$ws1.Range("N69").Formula = "=(L69-80000000)/80000000"
$ws1.Range("P69").Value = "(5932 with board and opposite caching)"

# ---------------------------------------------------------------------------
# 4. Move the viewport/selection down to the new bottom of the sheet.
# ---------------------------------------------------------------------------
$ws1.Range("A23").Select()
$ws1.Range("P69").Select()
